$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 377; this shifts the existing rows 377..472
# down to 378..473 (keeping all of their data intact), matching the diff's
# observed "cascade" of D/I/J/K/L/M/O/P values moving one row down.
$ws.Rows.Item(377).Insert()

# Populate the newly inserted row 377 with the new data record. Its static
# descriptive fields mirror what was already in the (now shifted-down) row
# 378, while the date (D) and volume (J) carry the new values from the diff.
$ws.Cells.Item(377, 1).Value = 9
$ws.Cells.Item(377, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(377, 3).Value = 'Metropolitana'
$ws.Cells.Item(377, 4).Value = 44782
$ws.Cells.Item(377, 5).Value = 13
$ws.Cells.Item(377, 6).Value = 100112012
$ws.Cells.Item(377, 7).Value = 'Espinaca'
$ws.Cells.Item(377, 8).Value = 'Sin especificar'
$ws.Cells.Item(377, 9).Value = 'Primera'
$ws.Cells.Item(377, 10).Value = 160
$ws.Cells.Item(377, 11).Value = 7000
$ws.Cells.Item(377, 12).Value = 8000
$ws.Cells.Item(377, 13).Value = 7500
$ws.Cells.Item(377, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(377, 15).Value = 'Provincia de Chacabuco'
$ws.Cells.Item(377, 16).Value = 750
$ws.Cells.Item(377, 17).Value = 10
$ws.Cells.Item(377, 18).Value = 'Hortaliza'
